$d = $word.ActiveDocument

foreach ($p in @($d.Paragraphs)) {
    if ($p.Range.Text -match "Sehr geehrte Damen und Herren") {
        $p.Range.Delete()
        break
    }
}
